# Trade #72 closed at 2026-02-18 00:27:29 - unknown UNKNOWN +0.000%
#
# Applies the trading-results update:
#   - Summary sheet roll-up numbers move to reflect the newly closed trade
#     and the newly opened trade.
#   - Strategy Status row for MarketMaking reflects the same roll-up.
#   - "All Trades" row 101 (Trade #100) flips from OPEN -> CLOSED (early_exit).
#   - "MarketMaking" row 33 (Trade #100) gets the same close-out.
#   - A brand-new OPEN trade (Trade #129) is appended as a new row on both
#     "All Trades" (row 130) and "MarketMaking" (row 50).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Summary
# ---------------------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B3").Value = 1499.25
$summary.Range("B4").Value = 0.36
$summary.Range("B6").Value = 100
$summary.Range("B8").Value = 37
$summary.Range("B9").Value = 47

# ---------------------------------------------------------------------------
# Strategy Status - MarketMaking row (row 6)
# ---------------------------------------------------------------------------
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("C6").Value = 99.40000000000001
$status.Range("D6").Value = 32
$status.Range("E6").Value = -0.41
$status.Range("F6").Value = -0.6
$status.Range("G6").Value = 46.88

# ---------------------------------------------------------------------------
# All Trades - close out Trade #100 (row 101)
# ---------------------------------------------------------------------------
$allTrades = $wb.Worksheets.Item("All Trades")
$allTrades.Cells.Item(101, 7).Value = 0.65
$allTrades.Cells.Item(101, 8).Value = "CLOSED"
$allTrades.Cells.Item(101, 9).Value = -1.5152
$allTrades.Cells.Item(101, 10).Value = -0.01
$allTrades.Cells.Item(101, 11).Value = 99.40000000000001
$allTrades.Cells.Item(101, 12).Value = "early_exit"
$allTrades.Cells.Item(101, 13).Value = 0.13

# New row 130 -> Trade #129 (OPEN)
# Row 129 is used as the donor for the two cells that need to stay literal
# text instead of being auto-coerced by the COM layer: the "yyyy-mm-dd"
# date string (would otherwise become a date serial) and the genuinely
# blank Exit Price / Exit Reason cells (an empty-string .Value assignment
# gets dropped entirely instead of persisting an empty cell).
$allTrades.Cells.Item(130, 1).Value = 129
$allTrades.Cells.Item(129, 2).Copy($allTrades.Cells.Item(130, 2))
$allTrades.Cells.Item(130, 3).Value = "00:27:23"
$allTrades.Cells.Item(130, 4).Value = "MarketMaking"
$allTrades.Cells.Item(130, 5).Value = "DOWN"
$allTrades.Cells.Item(130, 6).Value = 0.66
$allTrades.Cells.Item(129, 7).Copy($allTrades.Cells.Item(130, 7))
$allTrades.Cells.Item(130, 8).Value = "OPEN"
$allTrades.Cells.Item(130, 9).Value = 0
$allTrades.Cells.Item(130, 10).Value = 0
$allTrades.Cells.Item(130, 11).Value = 99.40967800952272
$allTrades.Cells.Item(129, 12).Copy($allTrades.Cells.Item(130, 12))
$allTrades.Cells.Item(130, 13).Value = 0
$allTrades.Cells.Item(130, 14).Value = 0
$allTrades.Cells.Item(130, 15).Value = 0
$allTrades.Cells.Item(130, 16).Value = 0.65
$allTrades.Cells.Item(130, 17).Value = "Wide spread capture: 392 bps vs avg 287 bps"

# ---------------------------------------------------------------------------
# MarketMaking strategy sheet - close out Trade #100 (row 33)
# ---------------------------------------------------------------------------
$mm = $wb.Worksheets.Item("MarketMaking")
$mm.Cells.Item(33, 7).Value = 0.65
$mm.Cells.Item(33, 8).Value = "CLOSED"
$mm.Cells.Item(33, 9).Value = -1.5152
$mm.Cells.Item(33, 10).Value = -0.01
$mm.Cells.Item(33, 11).Value = 99.40000000000001
$mm.Cells.Item(33, 16).Value = "early_exit"
$mm.Cells.Item(33, 17).Value = 0.13

# New row 50 -> Trade #129 (OPEN). Same donor-row trick as above, using
# row 49 (the prior last row) for the literal date string and the blank
# Exit Price / Exit Reason cells.
$mm.Cells.Item(50, 1).Value = 129
$mm.Cells.Item(49, 2).Copy($mm.Cells.Item(50, 2))
$mm.Cells.Item(50, 3).Value = "00:27:23"
$mm.Cells.Item(50, 4).Value = "MarketMaking"
$mm.Cells.Item(50, 5).Value = "DOWN"
$mm.Cells.Item(50, 6).Value = 0.66
$mm.Cells.Item(49, 7).Copy($mm.Cells.Item(50, 7))
$mm.Cells.Item(50, 8).Value = "OPEN"
$mm.Cells.Item(50, 9).Value = 0
$mm.Cells.Item(50, 10).Value = 0
$mm.Cells.Item(50, 11).Value = 99.40967800952272
$mm.Cells.Item(50, 12).Value = 0
$mm.Cells.Item(50, 13).Value = 0
$mm.Cells.Item(50, 14).Value = 0.65
$mm.Cells.Item(50, 15).Value = "Wide spread capture: 392 bps vs avg 287 bps"
$mm.Cells.Item(49, 16).Copy($mm.Cells.Item(50, 16))
$mm.Cells.Item(50, 17).Value = 0
